# Journal de bord PreTPI - apply "Template of the website + maquette + MCD-MLD review"
# Rows 3-13 get their task/date/time content rewritten in place; a new blank row
# is appended at the bottom of the table (row 80) pushing the closing footer row
# from 81 to 82; row35's style is normalised to match the surrounding band.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# 1) Fix up cell styles that need to change index (via format-only copy/paste
#    from donor cells whose style is not itself being altered).
# ---------------------------------------------------------------------------

# C5 must take on the "blank data row" look (same as C35/C36/.../C79)
$ws.Range("C35").Copy() | Out-Null
$ws.Range("C5").PasteSpecial($xlPasteFormats) | Out-Null

# D5 must take on the plain date style used by D7 (unchanged donor)
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D5").PasteSpecial($xlPasteFormats) | Out-Null

# C6 takes on the style that C4 keeps (style 19)
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C6").PasteSpecial($xlPasteFormats) | Out-Null

# D6 takes on the style that D4 keeps (style 11)
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D6").PasteSpecial($xlPasteFormats) | Out-Null

# C10 and C13 become "blank data row" look too
$ws.Range("C35").Copy() | Out-Null
$ws.Range("C10").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C35").Copy() | Out-Null
$ws.Range("C13").PasteSpecial($xlPasteFormats) | Out-Null

# ---------------------------------------------------------------------------
# 2) Update row heights: ht=30 moves from row 8 to row 9
# ---------------------------------------------------------------------------
$ws.Rows("8").AutoFit() | Out-Null
$ws.Rows("9").RowHeight = 30

# ---------------------------------------------------------------------------
# 3) Write the new cell text/values for rows 3-13
# ---------------------------------------------------------------------------

# Row 3 : only "Temps" changes (stays "45 minutes" but re-points in sst)
$ws.Range("E3").Value = "45 minutes"

# Row 4 : Task "MCD-MLD sur papier" -> "MCD sur papier"; Temps -> "20 minutes"
$ws.Range("C4").Value = "MCD sur papier"
$ws.Range("E4").Value = "20 minutes"

# Row 5 : brand-new row "MLD sur papier"
$ws.Range("B5").ClearContents() | Out-Null
$ws.Range("C5").Value = "MLD sur papier"
$ws.Range("D5").Value = 44228
$ws.Range("E5").Value = "100 minutes"

# Row 6 : old row-5 content ("Maison" / "Definition des sprints")
$ws.Range("B6").Value = "Maison"
$ws.Range("C6").Value = "Definition des sprints"
$ws.Range("D6").Value = 44229
$ws.Range("E6").Value = "45 minutes"

# Row 7 : "Création du journal de bord"
$ws.Range("C7").Value = "Création du journal de bord"
$ws.Range("D7").Value = 44231
$ws.Range("E7").Value = "10 minutes"
$ws.Range("F7").ClearContents() | Out-Null

# Row 8 : "Création de l'architecture" + description moved here
$ws.Range("C8").Value = "Création de l'architecture"
$ws.Range("D8").Value = 44231
$ws.Range("E8").Value = "10 minutes"
$ws.Range("F8").Value = "Dossier avec tout reuni + création de la docs"

# Row 9 : "MCD-MLD au format éléctronique + Correction"
$ws.Range("C9").Value = "MCD-MLD au format éléctronique + Correction"
$ws.Range("D9").Value = 44231
$ws.Range("E9").Value = "30 minutes"

# Row 10 : brand-new "Finition du sprint 1"
$ws.Range("C10").Value = "Finition du sprint 1"
$ws.Range("D10").Value = 44232
$ws.Range("E10").Value = "30 minutes"

# Row 11 : "Revision du sprint 1"
$ws.Range("C11").Value = "Revision du sprint 1"
$ws.Range("D11").Value = 44232
$ws.Range("E11").Value = "30 minutes"

# Row 12 : brand-new "Recherche du templates"
$ws.Range("C12").Value = "Recherche du templates"
$ws.Range("D12").Value = 44232
$ws.Range("E12").Value = "60 minutes"

# Row 13 : brand-new "Installation du template"
$ws.Range("C13").Value = "Installation du template"
$ws.Range("D13").Value = 44235
$ws.Range("E13").Value = "20 minutes"

# ---------------------------------------------------------------------------
# 4) Normalise row 35's style to the surrounding band (C35: s8 -> s20)
# ---------------------------------------------------------------------------
$ws.Range("C34").Copy() | Out-Null
$ws.Range("C35").PasteSpecial($xlPasteFormats) | Out-Null

# ---------------------------------------------------------------------------
# 5) Grow the table by one blank row at the bottom: old row 81 (footer)
#    becomes row 82; a new blank row 80 appears before it.
# ---------------------------------------------------------------------------
$ws.Range("C81:F81").Copy() | Out-Null
$ws.Range("C82:F82").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("C79:F79").Copy() | Out-Null
$ws.Range("C80:F80").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("C81:F81").ClearFormats() | Out-Null
$ws.Range("C81:F81").ClearContents() | Out-Null

# ---------------------------------------------------------------------------
# 6) Resize the table / autofilter to include the new row
# ---------------------------------------------------------------------------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("C2:F81"))

# ---------------------------------------------------------------------------
# 7) Selection moves to C14
# ---------------------------------------------------------------------------
$ws.Range("C14").Select() | Out-Null
